# Automatische test-sync: 2025-07-27 17:11:50
# Adds a new log row (row 5) to the "Logs" sheet, extends the conditional
# formatting ranges to cover it, and bumps the "Overig" tally on the
# "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 5 to the Logs sheet -----------------------------------
$logs.Range("A5").Value = "Kun jij dit even regelen?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D5").Value = "Overig"
$logs.Range("E5").Value = "Geachte klant,`nBedankt voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u wat meer informatie kunnen verstrekken over wat u precies wilt regelen?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F5").Value = "2025-07-27 17:11:00"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Ja"
$logs.Range("J5").Value = "Ja"

# Undo the automatic row-height bump that comes from the multi-line value in
# column E so row 5 stays on the default (un-customized) height like rows 2-4.
$logs.Rows.Item(5).AutoFit()

# --- Extend the conditional formatting ranges to include row 5 -----------
$ranges = @("D2:D4", "G2:G4", "H2:H4", "I2:I4", "J2:J4")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newAddr = "$col" + "2:" + "$col" + "5"
    $fcs = $logs.Range($addr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newAddr))
    }
}

# --- Update the Dashboard tally -------------------------------------------
$dashboard.Range("B2").Value = 4
